$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.Id -eq $id) { return $cand }
    }
    throw "Shape with id $id not found"
}

# Connector 60: shrink height (cy 927122 -> 545942)
$sh60 = Get-ShapeById $s 60
$sh60.Height = 42.98755905511811

# Shape 5: widen box to the left, add "Using tools" line
$sh5 = Get-ShapeById $s 5
$sh5.Left = 130.226066592126
$sh5.Width = 157.32944881889765
$tr5 = $sh5.TextFrame.TextRange
$tr5.Text = "Scanning Application or OS Image" + [char]13 + "Using tools: TerraScan, Regula, Checkov"
$p5_2 = $tr5.Paragraphs(2, 1)
$c = $p5_2.Characters(1, 13)
$c.Font.Color.RGB = 255
$c = $p5_2.Characters(14, 9)
$c.Font.Color.RGB = 255
$c = $p5_2.Characters(23, 10)
$c.Font.Color.RGB = 255
$c = $p5_2.Characters(33, 7)
$c.Font.Color.RGB = 255

# Connector 7: flip horizontally, reposition/resize
$sh7 = Get-ShapeById $s 7
$sh7.Left = 208.8907874015748
$sh7.Width = 27.168662077322836
$sh7.HorizontalFlip = [Microsoft.Office.Core.MsoTriState]::msoTrue

# Shape 15: move up + grow taller, add "Using tools" line
$sh15 = Get-ShapeById $s 15
$sh15.Top = 261.1809539818898
$sh15.Height = 107.08354568708661
$tr15 = $sh15.TextFrame.TextRange
$tr15.Text = "Scanning Infrastructure code change" + [char]13 + "Using tools: AWS Inspector, ECR Image Scanning, Trivy"
$p15_2 = $tr15.Paragraphs(2, 1)
$c = $p15_2.Characters(1, 48)
$c.Font.Color.RGB = 255
$c = $p15_2.Characters(49, 5)
$c.Font.Color.RGB = 255

# Connector 23: reposition/resize
$sh23 = Get-ShapeById $s 23
$sh23.Top = 314.7227630655118
$sh23.Height = 30.5103941007874

# Shape 28: grow box, add "Using tools:" + tool list lines
$sh28 = Get-ShapeById $s 28
$sh28.Width = 149.3844094488189
$sh28.Height = 73.7755928111811
$tr28 = $sh28.TextFrame.TextRange
$tr28.Text = "Post-Deployment Compliance Scan" + [char]13 + "Using tools:" + [char]13 + "AWS Config, AWS Security Hub"
$p28_2 = $tr28.Paragraphs(2, 1)
$p28_2.Font.Color.RGB = 255
$p28_3 = $tr28.Paragraphs(3, 1)
$p28_3.Font.Color.RGB = 255

# Connector 30: resize
$sh30 = Get-ShapeById $s 30
$sh30.Width = 23.196141732283465

